# Topic 6 -> Topic 8 (slide 1, shape "object 5" / Shapes.Item(3))
#
# The textbox currently reads "Topic<TAB>6" as two runs:
#   run1 = "T"          (rPr sz=3000 spc="-5" dirty="0")
#   run2 = "opic\t6"    (rPr sz=3000 dirty="0")
# We only need to turn the trailing "6" into "8"; editing just that
# character makes the host re-split run2 into "opic\t" + "8" (same rPr
# carried over to both pieces), exactly like PowerPoint does when you
# retype a single character in the middle of a run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Sanity check / locate the shape robustly instead of trusting a fixed index.
if ($tr.Text -notlike "*Topic*6*") {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $cand = $s.Shapes.Item($i)
        if ($cand.HasTextFrame -and $cand.TextFrame.TextRange.Text -like "*Topic*6*") {
            $sh = $cand
            $tr = $sh.TextFrame.TextRange
            break
        }
    }
}

$len = $tr.Length
$last = $tr.Characters($len, 1)
$last.Text = "8"

# The textbox auto-fits its height to the text (<a:spAutoFit/>); PowerPoint
# relays it out after the edit. Nudge past the float->EMU truncation so the
# stored height lands on the exact laid-out value instead of one EMU short.
$sh.Height = (461665 / 12700.0) + 0.00001
